# Edit script: applies the target diff to before.pptx
#  1. Slide 6's table gets a new (built-in) table style id.
#  2. The presentation's active theme colour scheme (currently the
#     "Integral" palette) is swapped for the "Office Theme" palette -
#     this is the colour-scheme portion of the theme1.xml / theme2.xml
#     content swap recorded in the diff.

$p = $ppt.ActivePresentation

# --- 1. Table style id on slide 6 -----------------------------------
$slide6 = $p.Slides.Item(6)
$tableShape = $slide6.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{5547C608-DC76-41DB-8392-A19AF3B74C6F}")

# --- 2. Swap the active theme's colour scheme (Integral -> Office) --
$scheme = $p.SlideMaster.ColorScheme
$scheme.Colors(1).RGB  = 0          # dk1       000000
$scheme.Colors(2).RGB  = 16777215   # lt1       FFFFFF
$scheme.Colors(3).RGB  = 6968388    # dk2       44546A
$scheme.Colors(4).RGB  = 15132391   # lt2       E7E6E6
$scheme.Colors(5).RGB  = 13998939   # accent1   5B9BD5
$scheme.Colors(6).RGB  = 3243501    # accent2   ED7D31
$scheme.Colors(7).RGB  = 10855845   # accent3   A5A5A5
$scheme.Colors(8).RGB  = 49407      # accent4   FFC000
$scheme.Colors(9).RGB  = 12874308   # accent5   4472C4
$scheme.Colors(10).RGB = 4697456    # accent6   70AD47
$scheme.Colors(11).RGB = 12673797   # hlink     0563C1
$scheme.Colors(12).RGB = 7491477    # folHlink  954F72
